$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value = "01-09-2021"
$ws.Range("A46").Style = "Normal"
$ws.Range("B46").Value = 111.45
$ws.Range("C46").Value = 109.43
$ws.Range("D46").Value = 113.37
$ws.Range("E46").Value = 109.18
$ws.Range("F46").Value = 122.49
